# Daily attendance processing - 2025-10-11 01:12:21
#
# The "Recorded By" column (G) stores a comma-separated list of the people/
# processes that recorded a given attendance session. This pass normalizes
# that list so that the automated "System" entry is always listed first
# (case preserved for each individual token), matching the convention used
# by the attendance pipeline. Rows whose "Recorded By" value does not
# contain a "system" token (any case) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in the "Recorded By" column G) that contain at least one "system"
# token and therefore need to be reordered.
$rows = @(2,3,4,5,6,7,10,11,12,13,14,15,29,30,32,33,34,37,38,39,40,41,42,56,57,58,59,60,61,64,65,66,67,68,69,84,85,86,87,88,89,93,95,110,111,112,113,114,115,119,121,136,137,138,139,140,141,145,147)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = [string]$cell.Value()
    $parts = $current -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    # Find every token that equals "system" (case-insensitive).
    $sysIndexes = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].ToLower() -eq 'system') {
            $sysIndexes += $i
        }
    }

    if ($sysIndexes.Length -ge 2) {
        # Two (or more) "system" tokens: swap the first and last occurrence,
        # keeping their original text/casing but exchanging position.
        $first = $sysIndexes[0]
        $last = $sysIndexes[$sysIndexes.Length - 1]
        $tmp = $parts[$first]
        $parts[$first] = $parts[$last]
        $parts[$last] = $tmp
        $newParts = $parts
    } elseif ($sysIndexes.Length -eq 1) {
        # Single "system" token: move it to the front of the list.
        $idx = $sysIndexes[0]
        $newParts = @($parts[$idx])
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $idx) {
                $newParts += $parts[$i]
            }
        }
    } else {
        $newParts = $parts
    }

    $newValue = [string]::Join(', ', $newParts)
    $cell.Value = $newValue
}
